$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new "longname" values for the Transport Canada Dash 7 / Dash 8 rows
$ws.Range("F5").Value = "Transport Canada Dash 7 - CGCFR"
$ws.Range("F6").Value = "Transport Canada Dash 8 - CGCFJ"

# New font/style for these two cells: Arial 12, color FF212529
$ws.Range("F5:F6").Font.Name = "Arial"
$ws.Range("F5:F6").Font.Size = 12
$ws.Range("F5:F6").Font.Color = 2696481

# Row heights for rows 5 and 6 changed to 15.75
$ws.Rows.Item(5).RowHeight = 15.75
$ws.Rows.Item(6).RowHeight = 15.75

# Update the selection to F6
$ws.Activate()
$ws.Range("F6").Select()
